# Apply "added few more programs 07/24" changes to practicePrograms.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New program file names added to columns B/C for rows that previously had
# blank Status/ProgramName cells. The shared-string table in the original
# commit has these new strings inserted in this exact order, so write the
# cell values in the same order to reproduce it.
$ws.Range("B8").Value = "swapWithout3Variable.java"
$ws.Range("B9").Value = "evenOrOdd.java"
$ws.Range("B10").Value = "vowelOrConsonant.java"
$ws.Range("B14").Value = "positiveOrNegative.java"
$ws.Range("B15").Value = "aplhabetOrNot.java"
$ws.Range("B16").Value = "sumOfNaturalNumbers.java"
$ws.Range("B11").Value = "largestAmongThree.java"

$ws.Range("C8").Value = "DONE"
$ws.Range("C9").Value = "DONE"
$ws.Range("C10").Value = "DONE"
$ws.Range("C11").Value = "DONE"
$ws.Range("C14").Value = "DONE"
$ws.Range("C15").Value = "DONE"
$ws.Range("C16").Value = "DONE"

# Update the view/selection state to match the edited workbook
# (new active cell is C11, with the view scrolled so row 7 is at the top).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("C11").Select()
